# Apply "repull data, push all data, mean calculation" update:
# Column F (dSF) values were repulled for a subset of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 1
    10 = -1
    12 = 4
    13 = 2
    14 = 5
    15 = 2
    18 = -1
    20 = 3
    28 = -3
    29 = 4
    32 = -2
    33 = -3
    42 = -5
    47 = -3
    50 = -4
    51 = -1
    52 = -5
    55 = -1
    56 = -2
    57 = -2
    62 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
